$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update date serial in A1
$ws.Range("A1").Value = 45309

# Step 1 & 2 price increases (~15%)
$ws.Range("D23").Value = 44.038
$ws.Range("D24").Value = 44.038
$ws.Range("D25").Value = 26.097
$ws.Range("D26").Value = 26.097
$ws.Range("D41").Value = 97.24299999999999
$ws.Range("D42").Value = 97.24299999999999
$ws.Range("D43").Value = 110.905
$ws.Range("D44").Value = 138.75
